$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2587.125
$ws.Range("J40").Value = 1857.2858
$ws.Range("L40").Value = 1857.2858
$ws.Range("N40").Value = -2207.2858
$ws.Range("H64").Value = 3638.2632
$ws.Range("I64").Value = 3385.3333
$ws.Range("J64").Value = 4071.8572
$ws.Range("K64").Value = 3385.3333
$ws.Range("L64").Value = 4071.8572
$ws.Range("M64").Value = -3137.3333
$ws.Range("N64").Value = -4567.8572
$ws.Range("H67").Value = 3638.2632
$ws.Range("I67").Value = 3385.3333
$ws.Range("J67").Value = 4071.8572
$ws.Range("K67").Value = 3385.3333
$ws.Range("L67").Value = 4071.8572
$ws.Range("M67").Value = -2527.3333
$ws.Range("N67").Value = -5787.8572
$ws.Range("H82").Value = 1925
$ws.Range("I82").Value = 850
$ws.Range("K82").Value = 2550
$ws.Range("M82").Value = -2144
$ws.Range("H85").Value = 1925
$ws.Range("I85").Value = 850
$ws.Range("K85").Value = 2550
$ws.Range("M85").Value = -1146
$ws.Range("H97").Value = 999.6
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 999.6
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2998.8
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -3990.8
$ws.Range("H100").Value = 1338.1818
$ws.Range("I100").Value = 1246.6666
$ws.Range("J100").Value = 1750
$ws.Range("K100").Value = 1246.6666
$ws.Range("L100").Value = 1750
$ws.Range("M100").Value = -705.6666
$ws.Range("N100").Value = -2832
$ws.Range("H118").Value = 195.85715
$ws.Range("I118").Value = 195.85715
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 587.5714499999999
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 1069.42855
$ws.Range("N118").ClearContents()
$ws.Range("H123").Value = 62000
$ws.Range("J123").Value = 62000
$ws.Range("L123").Value = 62000
$ws.Range("N123").Value = -71800
$ws.Range("H124").Value = 77400
$ws.Range("J124").Value = 77400
$ws.Range("L124").Value = 77400
$ws.Range("N124").Value = -87220
$ws.Range("H126").Value = 77900
$ws.Range("J126").Value = 77900
$ws.Range("L126").Value = 77900
$ws.Range("N126").Value = -87780
$ws.Range("H137").Value = 567650.5600000001
$ws.Range("I137").Value = 2527.4194
$ws.Range("J137").Value = 1193322.6
$ws.Range("K137").Value = 7582.2582
$ws.Range("L137").Value = 3579967.8
$ws.Range("M137").Value = -5032.2582
$ws.Range("N137").Value = -3585067.8
$ws.Range("H138").Value = 5159.8975
$ws.Range("I138").Value = 2093.8235
$ws.Range("J138").Value = 7529.136
$ws.Range("K138").Value = 6281.470499999999
$ws.Range("L138").Value = 22587.408
$ws.Range("M138").Value = -1141.470499999999
$ws.Range("N138").Value = -32867.408

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 50000
$ws.Range("J7").Value = 50000
$ws.Range("L7").Value = 50000
$ws.Range("N7").Value = -50228
$ws.Range("H32").Value = 21149.855
$ws.Range("I32").Value = 22764.56
$ws.Range("J32").Value = 5002.8
$ws.Range("K32").Value = 22764.56
$ws.Range("L32").Value = 5002.8
$ws.Range("M32").Value = -22477.56
$ws.Range("N32").Value = -5576.8
$ws.Range("H102").Value = 3800
$ws.Range("I102").Value = 3200
$ws.Range("J102").Value = 4400
$ws.Range("K102").Value = 3200
$ws.Range("L102").Value = 4400
$ws.Range("M102").Value = -1578
$ws.Range("N102").Value = -7644
$ws.Range("H135").Value = 41966.5
$ws.Range("J135").Value = 41966.5
$ws.Range("L135").Value = 41966.5
$ws.Range("N135").Value = -52106.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1003.2
$ws.Range("I99").Value = 919.0769
$ws.Range("J99").Value = 1550
$ws.Range("K99").Value = 919.0769
$ws.Range("L99").Value = 1550
$ws.Range("M99").Value = 578.9231
$ws.Range("N99").Value = -4546
$ws.Range("H103").Value = 16617.75
$ws.Range("J103").Value = 16617.75
$ws.Range("L103").Value = 16617.75
$ws.Range("N103").Value = -18961.75
$ws.Range("H105").Value = 8272.454
$ws.Range("J105").Value = 5142.857
$ws.Range("L105").Value = 5142.857
$ws.Range("N105").Value = -8636.857
$ws.Range("H107").Value = 1955.2222
$ws.Range("I107").Value = 1137.5
$ws.Range("J107").Value = 2609.4
$ws.Range("K107").Value = 1137.5
$ws.Range("L107").Value = 2609.4
$ws.Range("M107").Value = 782.5
$ws.Range("N107").Value = -6449.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 590772.9399999999
$ws.Range("I31").Value = 11256.526
$ws.Range("J31").Value = 852935.2
$ws.Range("K31").Value = 11256.526
$ws.Range("L31").Value = 852935.2
$ws.Range("M31").Value = -10961.526
$ws.Range("N31").Value = -853525.2
$ws.Range("H34").Value = 590772.9399999999
$ws.Range("I34").Value = 11256.526
$ws.Range("J34").Value = 852935.2
$ws.Range("K34").Value = 11256.526
$ws.Range("L34").Value = 852935.2
$ws.Range("M34").Value = -11054.526
$ws.Range("N34").Value = -853339.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2710.86
$ws.Range("I68").Value = 1519.2609
$ws.Range("J68").Value = 3725.926
$ws.Range("K68").Value = 4557.7827
$ws.Range("L68").Value = 11177.778
$ws.Range("M68").Value = -3746.7827
$ws.Range("N68").Value = -12799.778
$ws.Range("H71").Value = 2710.86
$ws.Range("I71").Value = 1519.2609
$ws.Range("J71").Value = 3725.926
$ws.Range("K71").Value = 13673.3481
$ws.Range("L71").Value = 33533.334
$ws.Range("M71").Value = -9617.348099999999
$ws.Range("N71").Value = -41645.334
$ws.Range("H99").Value = 2290.7144
$ws.Range("I99").Value = 1843.75
$ws.Range("J99").Value = 2886.6667
$ws.Range("K99").Value = 5531.25
$ws.Range("L99").Value = 8660.000100000001
$ws.Range("M99").Value = -3285.25
$ws.Range("N99").Value = -13152.0001
$ws.Range("H133").Value = 3640
$ws.Range("H134").Value = 4093.1365
$ws.Range("I134").Value = 4002.45
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 12007.35
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -6937.349999999999
$ws.Range("N134").Value = -25140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 79800
$ws.Range("J119").Value = 79800
$ws.Range("L119").Value = 79800
$ws.Range("N119").Value = -89476
$ws.Range("H132").Value = 6476.5
$ws.Range("I132").Value = 2946.4546
$ws.Range("J132").Value = 14242.6
$ws.Range("K132").Value = 8839.363799999999
$ws.Range("L132").Value = 42727.8
$ws.Range("M132").Value = -6309.363799999999
$ws.Range("N132").Value = -47787.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 727.7273
$ws.Range("I16").Value = 700.3333
$ws.Range("J16").Value = 851
$ws.Range("K16").Value = 700.3333
$ws.Range("L16").Value = 851
$ws.Range("M16").Value = -530.3333
$ws.Range("N16").Value = -1191
$ws.Range("H36").Value = 56000
$ws.Range("J36").Value = 56000
$ws.Range("L36").Value = 56000
$ws.Range("N36").Value = -57124
$ws.Range("H132").Value = 8005.3335
$ws.Range("I132").Value = 8886.053
$ws.Range("J132").Value = 4658.6
$ws.Range("K132").Value = 26658.159
$ws.Range("L132").Value = 13975.8
$ws.Range("M132").Value = -24128.159
$ws.Range("N132").Value = -19035.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 15923.667
$ws.Range("J74").Value = 16664.125
$ws.Range("L74").Value = 16664.125
$ws.Range("N74").Value = -18536.125
$ws.Range("H77").Value = 15923.667
$ws.Range("J77").Value = 16664.125
$ws.Range("L77").Value = 49992.375
$ws.Range("N77").Value = -59352.375
$ws.Range("H96").Value = 3374.5
$ws.Range("I96").Value = 1249
$ws.Range("J96").Value = 5500
$ws.Range("K96").Value = 1249
$ws.Range("L96").Value = 5500
$ws.Range("M96").Value = 124
$ws.Range("N96").Value = -8246
$ws.Range("H132").Value = 4250.523
$ws.Range("I132").Value = 4492.5137
$ws.Range("J132").Value = 2971.4285
$ws.Range("K132").Value = 13477.5411
$ws.Range("L132").Value = 8914.2855
$ws.Range("M132").Value = -10947.5411
$ws.Range("N132").Value = -13974.2855
